$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1): shift the displayed week from March to April.
#    B1:F1 hold date-like text (e.g. "2024.3.10"). Typing such a string
#    directly would make Excel auto-convert it to a serial date, so we:
#      a) stash the original (text) number format/style of B1:F1 on a
#         scratch range far outside the used area,
#      b) force the target range to Text format so the new strings are not
#         reinterpreted as dates,
#      c) write the new values,
#      d) paste the stashed formatting back on top (restores the original
#         style id), then wipe the scratch range.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "4月"

$ws.Range("B1:F1").Copy()
$ws.Range("H1:L1").PasteSpecial(-4122)

$ws.Range("B1:F1").NumberFormat = "@"
$ws.Range("B1").Value = "2024.4.7"
$ws.Range("C1").Value = "2024.4.8"
$ws.Range("D1").Value = "2024.4.9"
$ws.Range("E1").Value = "2024.4.10"
$ws.Range("F1").Value = "2024.4.11"

$ws.Range("H1:L1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("H1:L1").Clear()

# ---------------------------------------------------------------------------
# 2. Rows 2, 3, 5, 6, 7: the value that used to live in column D is moved to
#    column C (column D becomes blank), columns B/E/F are untouched.
# ---------------------------------------------------------------------------
foreach ($r in 2,3,5,6,7) {
    $ws.Cells.Item($r, 3).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 4).Value2 = ""
}

# ---------------------------------------------------------------------------
# 3. Row 4 (subjects) is rearranged: 数学/汉字/拼音 move one step around the
#    row, column D ends up blank.
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "数学"
$ws.Range("C4").Value = "汉字"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "拼音"
$ws.Range("F4").Value = "数学"

# ---------------------------------------------------------------------------
# 4. Update the selection shown in the sheet view.
# ---------------------------------------------------------------------------
$ws.Range("F2:F7").Select() | Out-Null
